$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings): edit only the changed run via Characters ---
$ws.Range("A8").Characters(21, 2).Text = "40"
$ws.Range("C9").Characters(27, 9).Text = "9/30/2024"
$ws.Range("C9").Characters(47, 9).Text = "10/6/2024"

# --- Cells that change from the text placeholder style (dash) to a numeric style: set value + number format ---
$ws.Range("D15").Value = 3
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C16").Value = 4
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 6
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("E16").Value = -33.333333333333
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D28").Value = 3
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -66.666666666666
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F31").Value = 4
$ws.Range("F31").NumberFormat = '#,##0'

# --- Cells that just change value, keeping their existing numeric style ---
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -27.777777777777
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 88
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = 27.536231884058
$ws.Range("L16").Value = 29.411764705882
$ws.Range("M16").Value = -37.142857142857
$ws.Range("N16").Value = -83.612662942271
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -81.818181818181
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -22.727272727272
$ws.Range("I17").Value = 164
$ws.Range("J17").Value = 169
$ws.Range("K17").Value = -2.958579881656
$ws.Range("L17").Value = -5.202312138728
$ws.Range("M17").Value = 37.815126050420
$ws.Range("N17").Value = -33.603238866396
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 62.5
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 2.173913043478
$ws.Range("L18").Value = -36.054421768707
$ws.Range("M18").Value = -67.918088737201
$ws.Range("N18").Value = -92.802450229709
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 419
$ws.Range("J19").Value = 437
$ws.Range("K19").Value = -4.118993135011
$ws.Range("L19").Value = -12.343096234309
$ws.Range("M19").Value = 42.517006802721
$ws.Range("N19").Value = -15.863453815261
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 47.058823529411
$ws.Range("I20").Value = 148
$ws.Range("J20").Value = 122
$ws.Range("K20").Value = 21.311475409836
$ws.Range("L20").Value = 64.444444444444
$ws.Range("M20").Value = 35.779816513761
$ws.Range("N20").Value = -89.367816091954
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -44.736842105263
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = 3.846153846153
$ws.Range("I21").Value = 930
$ws.Range("J21").Value = 909
$ws.Range("K21").Value = 2.310231023102
$ws.Range("L21").Value = -4.222451081359
$ws.Range("M21").Value = -3.326403326403
$ws.Range("N21").Value = -76.761619190404
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("M22").Value = -29.411764705882
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -5.555555555555
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 5.479452054794
$ws.Range("I24").Value = 751
$ws.Range("J24").Value = 833
$ws.Range("K24").Value = -9.843937575030
$ws.Range("L24").Value = -13.578826237054
$ws.Range("M24").Value = 5.477528089887
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -21.052631578947
$ws.Range("I25").Value = 185
$ws.Range("J25").Value = 246
$ws.Range("K25").Value = -24.796747967479
$ws.Range("L25").Value = -29.924242424242
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 330
$ws.Range("J26").Value = 279
$ws.Range("K26").Value = 18.279569892473
$ws.Range("L26").Value = 25.475285171102
$ws.Range("M26").Value = 0.303951367781
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = 15
$ws.Range("C28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = 1.923076923076
$ws.Range("L28").Value = -15.873015873015
$ws.Range("I31").Value = 15
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = 25

# --- Column E width change (bestFit width recalculated due to new wider values like "-100") ---
$ws.Columns("E:E").ColumnWidth = 7.433768
